$d = $word.ActiveDocument

# Locate the paragraph that contains the "(Guide: Polynomial long division)" text
# which immediately follows the "Further reading" heading, and remove the whole
# paragraph (including its paragraph mark).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "(Guide: Polynomial long division)") {
        $p.Range.Delete()
        break
    }
}

$d.Save()
